$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name reshuffle (shared-string table reorder) ---
$ws.Range("A117").Value = 'Angola'
$ws.Range("A118").Value = 'Cabo Verde'
$ws.Range("A119").Value = 'Lituania'
$ws.Range("A140").Value = 'Polinesia Francesa'
$ws.Range("A141").Value = 'Aruba'
$ws.Range("A142").Value = 'Mayotte'
$ws.Range("A143").Value = 'Islandia'
$ws.Range("A144").Value = 'Estonia'
$ws.Range("A145").Value = 'Somalia'

# --- Updated statistics ---
$ws.Range("B4").Value = 8447404
$ws.Range("C4").Value = 49303
$ws.Range("D4").Value = 5494518
$ws.Range("E4").Value = 2727715
$ws.Range("G4").Value = 397
$ws.Range("H4").Value = 225171
$ws.Range("B5").Value = 7594736
$ws.Range("C5").Value = 46498
$ws.Range("D5").Value = 6730617
$ws.Range("E5").Value = 748883
$ws.Range("G5").Value = 594
$ws.Range("H5").Value = 115236
$ws.Range("B6").Value = 5251127
$ws.Range("C6").Value = 15783
$ws.Range("E6").Value = 415242
$ws.Range("G6").Value = 321
$ws.Range("H6").Value = 154226
$ws.Range("B9").Value = 1002662
$ws.Range("C9").Value = 12982
$ws.Range("E9").Value = 171981
$ws.Range("G9").Value = 449
$ws.Range("H9").Value = 26716
$ws.Range("B12").Value = 870876
$ws.Range("C12").Value = 2201
$ws.Range("D12").Value = 784056
$ws.Range("E12").Value = 53000
$ws.Range("G12").Value = 61
$ws.Range("H12").Value = 33820
$ws.Range("B21").Value = 373731
$ws.Range("C21").Value = 6750
$ws.Range("E21").Value = 69032
$ws.Range("G21").Value = 33
$ws.Range("H21").Value = 9899
$ws.Range("D30").Value = 21158
$ws.Range("E30").Value = 190682
$ws.Range("B31").Value = 201437
$ws.Range("C31").Value = 3289
$ws.Range("D31").Value = 169671
$ws.Range("E31").Value = 21988
$ws.Range("G31").Value = 18
$ws.Range("H31").Value = 9778
$ws.Range("B34").Value = 181962
$ws.Range("C34").Value = 8077
$ws.Range("D34").Value = 74908
$ws.Range("E34").Value = 105541
$ws.Range("G34").Value = 91
$ws.Range("H34").Value = 1513
$ws.Range("B40").Value = 125181
$ws.Range("C40").Value = 436
$ws.Range("D40").Value = 101545
$ws.Range("E40").Value = 21062
$ws.Range("G40").Value = 10
$ws.Range("H40").Value = 2574
$ws.Range("B66").Value = 55452
$ws.Range("C66").Value = 728
$ws.Range("D66").Value = 36663
$ws.Range("E66").Value = 17582
$ws.Range("G66").Value = 19
$ws.Range("H66").Value = 1207
$ws.Range("B91").Value = 21506
$ws.Range("C91").Value = 65
$ws.Range("E91").Value = 965
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 424
$ws.Range("B101").Value = 13724
$ws.Range("C101").Value = 27
$ws.Range("E101").Value = 6124
$ws.Range("B113").Value = 8976
$ws.Range("C113").Value = 12
$ws.Range("D113").Value = 7303
$ws.Range("E113").Value = 1442
$ws.Range("B114").Value = 8884
$ws.Range("C114").Value = 3
$ws.Range("D114").Value = 8452
$ws.Range("E114").Value = 378
$ws.Range("B117").Value = 7829
$ws.Range("C117").Value = 207
$ws.Range("D117").Value = 3031
$ws.Range("E117").Value = 4550
$ws.Range("H117").Value = 248
$ws.Range("B118").Value = 7800
$ws.Range("C118").Value = 48
$ws.Range("D118").Value = 6620
$ws.Range("E118").Value = 1093
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 87
$ws.Range("B119").Value = 7726
$ws.Range("C119").Value = 205
$ws.Range("D119").Value = 3110
$ws.Range("E119").Value = 4503
$ws.Range("H119").Value = 113
$ws.Range("B120").Value = 7621
$ws.Range("C120").Value = 13
$ws.Range("D120").Value = 7355
$ws.Range("E120").Value = 103
$ws.Range("B134").Value = 5133
$ws.Range("C134").Value = 3
$ws.Range("D134").Value = 4959
$ws.Range("E134").Value = 65
$ws.Range("B140").Value = 4548
$ws.Range("C140").Value = 751
$ws.Range("D140").Value = 3202
$ws.Range("E140").Value = 1330
$ws.Range("G140").Value = 2
$ws.Range("H140").Value = 16
$ws.Range("B141").Value = 4334
$ws.Range("C141").Value = 12
$ws.Range("D141").Value = 4040
$ws.Range("E141").Value = 260
$ws.Range("H141").Value = 34
$ws.Range("B142").Value = 4159
$ws.Range("C142").Value = 129
$ws.Range("D142").Value = 2964
$ws.Range("E142").Value = 1152
$ws.Range("H142").Value = 43
$ws.Range("B143").Value = 4101
$ws.Range("C143").Value = 46
$ws.Range("D143").Value = 2856
$ws.Range("E143").Value = 1234
$ws.Range("H143").Value = 11
$ws.Range("B144").Value = 4085
$ws.Range("C144").Value = 7
$ws.Range("D144").Value = 3229
$ws.Range("E144").Value = 788
$ws.Range("H144").Value = 68
$ws.Range("B145").Value = 3890
$ws.Range("C145").Value = 26
$ws.Range("D145").Value = 3089
$ws.Range("E145").Value = 702
$ws.Range("H145").Value = 99
$ws.Range("B155").Value = 2560
$ws.Range("C155").Value = 29
$ws.Range("D155").Value = 2121
$ws.Range("E155").Value = 388
$ws.Range("B157").Value = 2403
$ws.Range("C157").Value = 14
$ws.Range("D157").Value = 1818
$ws.Range("E157").Value = 544
$ws.Range("B184").Value = 419
$ws.Range("C184").Value = 2
$ws.Range("D184").Value = 379
$ws.Range("E184").Value = 30

# --- Updated timestamp title ---
$ws.Range("A1").Value = 'Datos actualizados a 20 de Octubre de 2020 a las 01:57'
